# Re-label the "ScreenRecStarted" state as "0_unstated" throughout the
# transition-matrix sheet, and move the active selection to G1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header (row 1) that used to read "ScreenRecStarted"
$ws.Range("G1").Value = "0_unstated"

# Row labels (column A) that referenced "ScreenRecStarted"
$ws.Range("A27").Value = "0_unstated1_Scanning"
$ws.Range("A28").Value = "0_unstated3_Reading"
$ws.Range("A29").Value = "0_unstated5_Unknown "
$ws.Range("A30").Value = "0_unstated0_unstated"

# Update the active cell/selection shown when the sheet is opened.
$ws.Range("G1").Select()
